$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 259: drop the trailing " (ตัง)" from the question text
$ws.Cells.Item(259, 1).Value2 = "ถูกหนึ่ง"
$ws.Cells.Item(259, 2).Value2 = "ประเทศไทย เพราะประเทศไทยมีตรัง"

# Insert new training examples for tag "ถูกหนึ่ง" in rows 260-262
$ws.Cells.Item(260, 1).Value2 = "ถูกหนึ่ง"
$ws.Cells.Item(260, 2).Value2 = "ไทย เพราะประเทศไทยมีตรัง"

$ws.Cells.Item(261, 1).Value2 = "ถูกหนึ่ง"
$ws.Cells.Item(261, 2).Value2 = "ไทย"

$ws.Cells.Item(262, 1).Value2 = "ถูกหนึ่ง"
$ws.Cells.Item(262, 2).Value2 = "ตรัง"

# Row 263 intentionally left blank

# Move the former row 260 ("ถูกสอง" / "ไปฉันเพล") down to row 264
$ws.Cells.Item(264, 1).Value2 = "ถูกสอง"
$ws.Cells.Item(264, 2).Value2 = "ไปฉันเพล"

# Match the final cursor/selection position recorded in the saved workbook
[void]$ws.Range("C269").Select()

